# Fruta / hortaliza, semanal
# Insert a new weekly price-report row at row 296 (pushing the existing
# rows 296:315 down to 297:316) on the single data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 296:315 down to 297:316, leaving a blank row 296.
$ws.Rows.Item(296).Insert()

# Populate the newly inserted row 296 with the new record.
$ws.Range("A296").Value = 4
$ws.Range("B296").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C296").Value = "Los Lagos"
$ws.Range("D296").Value = 44746
$ws.Range("E296").Value = 10
$ws.Range("F296").Value = "Fruta"
$ws.Range("G296").Value = 100102
$ws.Range("H296").Value = "Cítricos"
$ws.Range("I296").Value = 100102006
$ws.Range("J296").Value = "Pomelo"
$ws.Range("K296").Value = "Start Ruby"
$ws.Range("L296").Value = "Primera"
$ws.Range("M296").Value = 60
$ws.Range("N296").Value = 13000
$ws.Range("O296").Value = 14000
$ws.Range("P296").Value = 13500
$ws.Range("Q296").Value = "$/caja 14 kilos empedrada"
$ws.Range("R296").Value = "Región de O'Higgins"
$ws.Range("S296").Value = 964
$ws.Range("T296").Value = 14
